$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 70.5
$ws.Range("I2").Value = 70.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 70.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 42.5
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 210.375
$ws.Range("I33").Value = 116.46154
$ws.Range("J33").Value = 617.3333
$ws.Range("K33").Value = 116.46154
$ws.Range("L33").Value = 617.3333
$ws.Range("M33").Value = 112.53846
$ws.Range("N33").Value = -1075.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 161707.78
$ws.Range("I129").Value = 801323.6
$ws.Range("J129").Value = 1803.8334
$ws.Range("K129").Value = 2403970.8
$ws.Range("L129").Value = 5411.5002
$ws.Range("M129").Value = -2398970.8
$ws.Range("N129").Value = -15411.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 21701.404
$ws.Range("I132").Value = 3547.5952
$ws.Range("J132").Value = 174193.4
$ws.Range("K132").Value = 10642.7856
$ws.Range("L132").Value = 522580.2
$ws.Range("M132").Value = -8112.785600000001
$ws.Range("N132").Value = -527640.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26911.395
$ws.Range("I32").Value = 26482.408
$ws.Range("J32").Value = 32574
$ws.Range("K32").Value = 26482.408
$ws.Range("L32").Value = 32574
$ws.Range("M32").Value = -26195.408
$ws.Range("N32").Value = -33148

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 13525
$ws.Range("J92").Value = 13525
$ws.Range("L92").Value = 13525
$ws.Range("N92").Value = -18517

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 45643.637
$ws.Range("J101").Value = 45643.637
$ws.Range("L101").Value = 45643.637
$ws.Range("N101").Value = -52133.637

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 350016670
$ws.Range("J112").Value = 350016670
$ws.Range("L112").Value = 350016670
$ws.Range("N112").Value = -350019624

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2751
$ws.Range("I122").Value = 2209.2222
$ws.Range("J122").Value = 3970
$ws.Range("K122").Value = 6627.6666
$ws.Range("L122").Value = 11910
$ws.Range("M122").Value = -4177.6666
$ws.Range("N122").Value = -16810

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2237.0435
$ws.Range("I20").Value = 2275.5
$ws.Range("K20").Value = 2275.5
$ws.Range("M20").Value = -2028.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 5625.5
$ws.Range("I22").Value = 10000
$ws.Range("J22").Value = 3000.8
$ws.Range("K22").Value = 10000
$ws.Range("L22").Value = 3000.8
$ws.Range("M22").Value = -9827
$ws.Range("N22").Value = -3346.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 34300
$ws.Range("J100").Value = 34300
$ws.Range("L100").Value = 34300
$ws.Range("N100").Value = -36464

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2033.3636
$ws.Range("I107").Value = 1827.421
$ws.Range("J107").Value = 3337.6667
$ws.Range("K107").Value = 1827.421
$ws.Range("L107").Value = 3337.6667
$ws.Range("M107").Value = 92.57899999999995
$ws.Range("N107").Value = -7177.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 8000
$ws.Range("I113").Value = 8000
$ws.Range("K113").Value = 8000
$ws.Range("M113").Value = -5830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 44985.6
$ws.Range("J133").Value = 44985.6
$ws.Range("L133").Value = 44985.6
$ws.Range("N133").Value = -55105.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2623.5
$ws.Range("I134").Value = 1611.2188
$ws.Range("J134").Value = 4423.1113
$ws.Range("K134").Value = 4833.6564
$ws.Range("L134").Value = 13269.3339
$ws.Range("M134").Value = -2298.6564
$ws.Range("N134").Value = -18339.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5442.7544
$ws.Range("I31").Value = 2490.3845
$ws.Range("J31").Value = 6315.0454
$ws.Range("K31").Value = 2490.3845
$ws.Range("L31").Value = 6315.0454
$ws.Range("M31").Value = -2195.3845
$ws.Range("N31").Value = -6905.0454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5442.7544
$ws.Range("I34").Value = 2490.3845
$ws.Range("J34").Value = 6315.0454
$ws.Range("K34").Value = 2490.3845
$ws.Range("L34").Value = 6315.0454
$ws.Range("M34").Value = -2288.3845
$ws.Range("N34").Value = -6719.0454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 46597
$ws.Range("J92").Value = 46597
$ws.Range("L92").Value = 46597
$ws.Range("N92").Value = -51589

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 72996.8
$ws.Range("J96").Value = 72996.8
$ws.Range("L96").Value = 72996.8
$ws.Range("N96").Value = -78488.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 25017
$ws.Range("J106").Value = 25017
$ws.Range("L106").Value = 25017
$ws.Range("N106").Value = -27541

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 47618.773
$ws.Range("I132").Value = 1713.7142
$ws.Range("K132").Value = 5141.142599999999
$ws.Range("M132").Value = -2611.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 12810
$ws.Range("J133").Value = 12810
$ws.Range("L133").Value = 12810
$ws.Range("N133").Value = -17870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3886.3333
$ws.Range("I134").Value = 1683
$ws.Range("J134").Value = 5355.222
$ws.Range("K134").Value = 5049
$ws.Range("L134").Value = 16065.666
$ws.Range("M134").Value = -2514
$ws.Range("N134").Value = -21135.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 901.3333
$ws.Range("J55").Value = 1500
$ws.Range("L55").Value = 4500
$ws.Range("N55").Value = -4854

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 7008.92
$ws.Range("J88").Value = 7008.92
$ws.Range("L88").Value = 21026.76
$ws.Range("N88").Value = -21882.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H91").Value = 7008.92
$ws.Range("J91").Value = 7008.92
$ws.Range("L91").Value = 21026.76
$ws.Range("N91").Value = -23990.76

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2601.5151
$ws.Range("I131").Value = 17134.666
$ws.Range("J131").Value = 1148.2
$ws.Range("K131").Value = 51403.99800000001
$ws.Range("L131").Value = 3444.6
$ws.Range("M131").Value = -46363.99800000001
$ws.Range("N131").Value = -13524.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5382.143
$ws.Range("I70").Value = 5694.706
$ws.Range("J70").Value = 4899.091
$ws.Range("K70").Value = 5694.706
$ws.Range("L70").Value = 4899.091
$ws.Range("M70").Value = -5424.706
$ws.Range("N70").Value = -5439.091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5382.143
$ws.Range("I73").Value = 5694.706
$ws.Range("J73").Value = 4899.091
$ws.Range("K73").Value = 5694.706
$ws.Range("L73").Value = 4899.091
$ws.Range("M73").Value = -4758.706
$ws.Range("N73").Value = -6771.091

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 11964.1
$ws.Range("I97").Value = 2500
$ws.Range("J97").Value = 13015.667
$ws.Range("K97").Value = 2500
$ws.Range("L97").Value = 13015.667
$ws.Range("M97").Value = -2004
$ws.Range("N97").Value = -14007.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2292.0732
$ws.Range("I132").Value = 1816.7354
$ws.Range("K132").Value = 5450.206200000001
$ws.Range("M132").Value = -2920.206200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1231.9333
$ws.Range("I22").Value = 1331.1111
$ws.Range("J22").Value = 1083.1666
$ws.Range("K22").Value = 1331.1111
$ws.Range("L22").Value = 1083.1666
$ws.Range("M22").Value = -1036.1111
$ws.Range("N22").Value = -1673.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1231.9333
$ws.Range("I27").Value = 1331.1111
$ws.Range("J27").Value = 1083.1666
$ws.Range("K27").Value = 1331.1111
$ws.Range("L27").Value = 1083.1666
$ws.Range("M27").Value = -1224.1111
$ws.Range("N27").Value = -1297.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 29000
$ws.Range("J41").Value = 29000
$ws.Range("L41").Value = 29000
$ws.Range("N41").Value = -29876

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2072.5789
$ws.Range("I46").Value = 1099.875
$ws.Range("J46").Value = 2780
$ws.Range("K46").Value = 1099.875
$ws.Range("L46").Value = 2780
$ws.Range("M46").Value = -911.875
$ws.Range("N46").Value = -3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 33736.168
$ws.Range("J104").Value = 33736.168
$ws.Range("L104").Value = 33736.168
$ws.Range("N104").Value = -40724.168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 50670.4
$ws.Range("J116").Value = 50670.4
$ws.Range("L116").Value = 50670.4
$ws.Range("N116").Value = -59848.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 43030.4
$ws.Range("J117").Value = 43030.4
$ws.Range("L117").Value = 43030.4
$ws.Range("N117").Value = -52208.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3807.639
$ws.Range("I132").Value = 3455.037
$ws.Range("J132").Value = 4865.4443
$ws.Range("K132").Value = 10365.111
$ws.Range("L132").Value = 14596.3329
$ws.Range("M132").Value = -7835.110999999999
$ws.Range("N132").Value = -19656.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 31333.334
$ws.Range("J101").Value = 31333.334
$ws.Range("L101").Value = 31333.334
$ws.Range("N101").Value = -37823.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 47328
$ws.Range("J110").Value = 47328
$ws.Range("L110").Value = 47328
$ws.Range("N110").Value = -55508
